$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column B (shifts old B..F to C..G), carrying A's date style along the row ---
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = 11
$ws.Columns.Item(6).ColumnWidth = 96.8

# --- Header row ---
$ws.Range("B1").Value = "from"
$ws.Range("H1").Value = "problem_at_which_end"

# --- New column H values for existing rows 2-9 ---
$ws.Range("H2").Value = "ISI"
$ws.Range("H3").Value = "ISI"
$ws.Range("H4").Value = "ISI"
$ws.Range("H5").Value = "PH"
$ws.Range("H6").Value = "none"
$ws.Range("H7").Value = "PH"
$ws.Range("H8").Value = "ISI"
$ws.Range("H9").Value = "ISI"

# --- New rows 11-15 ---

# Row 11
$ws.Range("A9").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 43199
$ws.Range("B11").Value = "Amandeep"
$ws.Range("C11").Value = 1028
$ws.Range("F11").Value = "5 answer options, but only 4 get back; related to monotonicity check in separations I added"
$ws.Range("H11").Value = "PH"

# Row 12
$ws.Range("A9").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 43201
$ws.Range("B12").Value = "Amandeep"
$ws.Range("C12").Value = 839
$ws.Range("E12").Value = "Will ACLED record any riot/protest events in Gabon in April 2018? "
$ws.Range("F12").Value = "R error; parse_separations monotonic check choked on length 1 input for binary IFP"
$ws.Range("H12").Value = "PH"

# Row 13
$ws.Range("A9").Copy($ws.Range("A13"))
$ws.Range("A13").Value = 43202
$ws.Range("B13").Value = "Amandeep"
$ws.Range("C13").Value = 1235
$ws.Range("E13").Value = "What will be the monthly period-over-period change in the consumer price index (CPI) for Benin in April 2018?"
$ws.Range("F13").Value = "R error; parse_separations did not recognize negative cutpoints"
$ws.Range("H13").Value = "PH"

# Row 14
$ws.Range("A9").Copy($ws.Range("A14"))
$ws.Range("A14").Value = 43202
$ws.Range("B14").Value = "Amandeep"
$ws.Range("C14").Value = 1406
$ws.Range("E14").Value = "What will be the monthly period-over-period change in the consumer price index (CPI) for Egypt in May 2018?"
$ws.Range("F14").Value = "R error; parse_separations did not recognize negative cutpoints"
$ws.Range("H14").Value = "PH"

# Row 15
$ws.Range("A9").Copy($ws.Range("A15"))
$ws.Range("A15").Value = 43202
$ws.Range("A9").Copy($ws.Range("B15"))
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = 1226
$ws.Range("E15").Value = "Will ACLED record any civilian fatalities in Ghana in June 2018?"
$ws.Range("F15").Value = "CI was in negative values; series not recognized as count because ACLED keyword was upper case, but ifp name converted to lower when matching"
$ws.Range("H15").Value = "PH"

# --- Selection matches the author's final cursor position ---
[void]$ws.Range("H13").Select()
